# Checklist_Prompts_Sistemas.xlsx edit
# "Audios sistemas generativos recortados"
# Replace the "Disco" genre block (rows 7-11) with a "Jazz" genre block on all
# three sheets (MUSICGEN, AIVA, SUNO), marking those rows as completed with
# their corresponding number of generated versions / observations, matching
# the pattern already used by the other genres (Classical, Metal) in each sheet.

$wb = $excel.ActiveWorkbook

# New Jazz prompts (order matters for shared-string bookkeeping: write the
# prompts first, then the genre label, then the last prompt, mirroring how
# the rest of the row was authored).
$jazzPrompt1 = "Generate a 30-second jazz music song."
$jazzPrompt2 = "Generate a 30-second jazz instrumental piece."
$jazzPrompt3 = "Generate a 30-second bebop jazz piece inspired by Charlie Parker."
$jazzPrompt4 = "Generate a 30-second swing jazz song in the style of Duke Ellington."
$jazzPrompt5 = "Generate a 30-second jazz song that sounds melancholic and nostalgic."
$check = "✅"
$seed = "Se usan estilos predeterminados parecidos como semilla"

function Update-GenreBlock {
    param([string]$SheetName, [int]$Versions, [bool]$FillObservations)

    $ws = $wb.Worksheets.Item($SheetName)

    # Prompts (column C) first
    $ws.Range("C7").Value = $jazzPrompt1
    $ws.Range("C8").Value = $jazzPrompt2
    $ws.Range("C9").Value = $jazzPrompt3
    $ws.Range("C10").Value = $jazzPrompt4

    # Genre label (column A) -- Disco -> Jazz
    $ws.Range("A7:A11").Value = "Jazz"

    # Last prompt written after the genre label, matching original authoring order
    $ws.Range("C11").Value = $jazzPrompt5

    # Mark rows as completed + number of versions generated (previously blank
    # for the now-removed Disco rows).
    $ws.Range("D7:D11").Value = $check
    $ws.Range("D7:D11").HorizontalAlignment = -4108   # xlCenter, matches sibling genre rows

    $ws.Range("E7").Value = $Versions
    $ws.Range("E8").Value = $Versions
    $ws.Range("E9").Value = $Versions
    $ws.Range("E10").Value = $Versions
    $ws.Range("E11").Value = $Versions

    if ($FillObservations) {
        $ws.Range("F7").Value = $seed
        $ws.Range("F8").Value = $seed
        $ws.Range("F9").Value = $seed
        $ws.Range("F10").Value = $seed
        $ws.Range("F11").Value = $seed
    }
}

Update-GenreBlock "MUSICGEN" 1 $false
Update-GenreBlock "AIVA"     1 $true
Update-GenreBlock "SUNO"     2 $false

# Restore the last-selected cell on each sheet to match the saved state.
$wsM = $wb.Worksheets.Item("MUSICGEN")
$wsM.Range("I8").Select()

$wsA = $wb.Worksheets.Item("AIVA")
$wsA.Range("C18").Select()

$wsS = $wb.Worksheets.Item("SUNO")
$wsS.Range("E18").Select()

$wsA.Activate()
